$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-43: price (D) and volume (E) updates ---
$ws.Range("D2").Value = "22.393.92"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.566.28"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9998"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.68"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3745"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.43"
$ws.Range("E9").Value = "  -5.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.150"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07411"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.46"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.826"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "1.564.93"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.65"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.352"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.25"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.68"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").Value = "22.402.10"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.304"
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.544"
$ws.Range("E26").Value = "  -2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.59"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.35"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.909"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.35"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "1.740.00"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.939"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.915"
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.621"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08254"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.298"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06311"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.247"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.09"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6088"
$ws.Range("E43").Value = "  -2.81%  "

# --- Rows 44-51: new Frax row inserted, content shifted down, Aave dropped ---
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.68"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.744"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5908"
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.008"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.60"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.180"
$ws.Range("E50").Value = "  -3.57%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07142"
$ws.Range("E51").Value = "  -1.09%  "
